$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item(1)

# Remove existing hyperlinks so stale relationship targets are not left behind
$ws.Hyperlinks.Delete()

# Row 2: 大手SIer等のAIソリューション開発・導入を支援してくださ
$ws.Range("A2").Value = '2025-12-17 01:20:50'
$ws.Range("B2").Value = '大手SIer等のAIソリューション開発・導入を支援してくださるエンジニア・PM募集'
$ws.Range("C2").Value = 'システム開発'
$ws.Range("D2").Value = '300,000 円 ~ 500,000 円 / 固定'
$ws.Range("E2").Value = '期限情報なし'
$ws.Range("F2").Value = 'https://www.lancers.jp/work/detail/5455098'
$ws.Range("G2").Value = 375
$ws.Range("H2").Value = '🔥AI,Ai ◆開発'

# Row 3: 【フルリモート】官公庁向けPythonアプリ開発PM募集|7
$ws.Range("A3").Value = '2025-12-17 01:20:50'
$ws.Range("B3").Value = '【フルリモート】官公庁向けPythonアプリ開発PM募集|7名チーム統括'
$ws.Range("C3").Value = 'システム開発'
$ws.Range("D3").Value = '500,000 円 ~ 1,000,000 円 / 固定'
$ws.Range("E3").Value = '期限情報なし'
$ws.Range("F3").Value = 'https://www.lancers.jp/work/detail/5454985'
$ws.Range("G3").Value = 295
$ws.Range("H3").Value = '🔥Python ◆開発 ◇アプリ'

# Row 4: 【せどり×ツール製作】APIを使用したせどりツールを製作でき
$ws.Range("A4").Value = '2025-12-17 01:20:50'
$ws.Range("B4").Value = '【せどり×ツール製作】APIを使用したせどりツールを製作できるエンジニアさんを募集します♪'
$ws.Range("C4").Value = 'システム開発'
$ws.Range("D4").Value = '20,000 円 ~ 50,000 円 / 固定'
$ws.Range("E4").Value = '期限情報なし'
$ws.Range("F4").Value = 'https://www.lancers.jp/work/detail/5217096'
$ws.Range("G4").Value = 243
$ws.Range("H4").Value = '🔥API ◆ツール'

# Row 5: 【Java/対話システム/心理学実験】協同問題解決プラットフ
$ws.Range("A5").Value = '2025-12-17 01:20:50'
$ws.Range("B5").Value = '【Java/対話システム/心理学実験】協同問題解決プラットフォームの改修開発'
$ws.Range("C5").Value = 'システム開発'
$ws.Range("D5").Value = '500,000 円 ~ 1,000,000 円 / 固定'
$ws.Range("E5").Value = '期限情報なし'
$ws.Range("F5").Value = 'https://www.lancers.jp/work/detail/5439921'
$ws.Range("G5").Value = 155
$ws.Range("H5").Value = '★Java ◆開発'

# Row 6: ホットペッパービューティーブログ一括投稿システム開発
$ws.Range("A6").Value = '2025-12-17 01:20:50'
$ws.Range("B6").Value = 'ホットペッパービューティーブログ一括投稿システム開発'
$ws.Range("C6").Value = 'システム開発'
$ws.Range("D6").Value = '20,000 円 ~ 50,000 円 / 固定'
$ws.Range("E6").Value = '期限情報なし'
$ws.Range("F6").Value = 'https://www.lancers.jp/work/detail/5455160'
$ws.Range("G6").Value = 113
$ws.Range("H6").Value = '◆開発,システム開発'

# Row 7: 【急募】新規システム開発に伴う要件定義依頼
$ws.Range("A7").Value = '2025-12-17 01:20:50'
$ws.Range("B7").Value = '【急募】新規システム開発に伴う要件定義依頼'
$ws.Range("C7").Value = 'システム開発'
$ws.Range("D7").Value = '10,000 円 ~ 20,000 円 / 固定'
$ws.Range("E7").Value = '期限情報なし'
$ws.Range("F7").Value = 'https://www.lancers.jp/work/detail/5455415'
$ws.Range("G7").Value = 110
$ws.Range("H7").Value = '◆開発,システム開発'

# Row 8: 【急募】iPhone・Android対応の天気アプリ開発をお
$ws.Range("A8").Value = '2025-12-17 01:20:50'
$ws.Range("B8").Value = '【急募】iPhone・Android対応の天気アプリ開発をお願いします!'
$ws.Range("C8").Value = 'システム開発'
$ws.Range("D8").Value = '500,000 円 ~ 1,000,000 円 / 固定'
$ws.Range("E8").Value = '期限情報なし'
$ws.Range("F8").Value = 'https://www.lancers.jp/work/detail/5455038'
$ws.Range("G8").Value = 100
$ws.Range("H8").Value = '◆開発 ◇アプリ'

# Row 9: Kabuステーション自動売買アプリの開発
$ws.Range("A9").Value = '2025-12-17 01:20:50'
$ws.Range("B9").Value = 'Kabuステーション自動売買アプリの開発'
$ws.Range("C9").Value = 'システム開発'
$ws.Range("D9").Value = '50,000 円 ~ 100,000 円 / 固定'
$ws.Range("E9").Value = '期限情報なし'
$ws.Range("F9").Value = 'https://www.lancers.jp/work/detail/5455251'
$ws.Range("G9").Value = 93
$ws.Range("H9").Value = '◆開発 ◇アプリ'

# Row 10: ホームページ診断チェックツール
$ws.Range("A10").Value = '2025-12-17 01:20:50'
$ws.Range("B10").Value = 'ホームページ診断チェックツール'
$ws.Range("C10").Value = 'システム開発'
$ws.Range("D10").Value = '50,000 円 ~ 100,000 円 / 固定'
$ws.Range("E10").Value = '期限情報なし'
$ws.Range("F10").Value = 'https://www.lancers.jp/work/detail/5455029'
$ws.Range("G10").Value = 73
$ws.Range("H10").Value = '◆ツール'

# Row 11: 【急募】帳票デジタル化のフロントエンド開発者募集
$ws.Range("A11").Value = '2025-12-17 01:20:50'
$ws.Range("B11").Value = '【急募】帳票デジタル化のフロントエンド開発者募集'
$ws.Range("C11").Value = 'システム開発'
$ws.Range("D11").Value = '50,000 円 ~ 100,000 円 / 固定'
$ws.Range("E11").Value = '期限情報なし'
$ws.Range("F11").Value = 'https://www.lancers.jp/work/detail/5454857'
$ws.Range("G11").Value = 68
$ws.Range("H11").Value = '◆開発'

# Row 12: 【急募】Accessシステム改修・CSV読込・MySQLクラ
$ws.Range("A12").Value = '2025-12-17 01:20:50'
$ws.Range("B12").Value = '【急募】Accessシステム改修・CSV読込・MySQLクラウド化・PDFデータ調整'
$ws.Range("C12").Value = 'システム開発'
$ws.Range("D12").Value = '50,000 円 ~ 100,000 円 / 固定'
$ws.Range("E12").Value = '期限情報なし'
$ws.Range("F12").Value = 'https://www.lancers.jp/work/detail/5455015'
$ws.Range("G12").Value = 53
$ws.Range("H12").Value = '◇MySQL'

# Row 13: 【改善提案募集】事業管理スプレッドシートの見直し・改善提案を
$ws.Range("A13").Value = '2025-12-17 01:20:50'
$ws.Range("B13").Value = '【改善提案募集】事業管理スプレッドシートの見直し・改善提案をお願いします。'
$ws.Range("C13").Value = 'システム開発'
$ws.Range("D13").Value = '1,000 ~ 5,000 円 / 固定'
$ws.Range("E13").Value = '期限情報なし'
$ws.Range("F13").Value = 'https://www.lancers.jp/work/detail/5455422'
$ws.Range("G13").Value = 30
$ws.Range("H13").Value = '◇管理'

# Row 14: 【急募】wixシステムでのメッセージ送信システム構築依頼
$ws.Range("A14").Value = '2025-12-17 01:20:50'
$ws.Range("B14").Value = '【急募】wixシステムでのメッセージ送信システム構築依頼'
$ws.Range("C14").Value = 'システム開発'
$ws.Range("D14").Value = '50,000 円 ~ 100,000 円 / 固定'
$ws.Range("E14").Value = '期限情報なし'
$ws.Range("F14").Value = 'https://www.lancers.jp/work/detail/5455067'
$ws.Range("G14").Value = 33
$ws.Range("H14").Value = ""

# Row 15: 【急募】企業のセキュリティ対策を担うエンジニア募集
$ws.Range("A15").Value = '2025-12-17 01:20:50'
$ws.Range("B15").Value = '【急募】企業のセキュリティ対策を担うエンジニア募集'
$ws.Range("C15").Value = 'システム開発'
$ws.Range("D15").Value = '500,000 円 ~ 1,000,000 円 / 固定'
$ws.Range("E15").Value = '期限情報なし'
$ws.Range("F15").Value = 'https://www.lancers.jp/work/detail/5450345'
$ws.Range("G15").Value = 25
$ws.Range("H15").Value = ""

# Row 16: 〖リモート可〗Delphiエンジニア募集
$ws.Range("A16").Value = '2025-12-17 01:20:50'
$ws.Range("B16").Value = '〖リモート可〗Delphiエンジニア募集'
$ws.Range("C16").Value = 'システム開発'
$ws.Range("D16").Value = '300,000 円 ~ 500,000 円 / 固定'
$ws.Range("E16").Value = '期限情報なし'
$ws.Range("F16").Value = 'https://www.lancers.jp/work/detail/5341051'
$ws.Range("G16").Value = 25
$ws.Range("H16").Value = ""

# Row 17: 【SESエンジニア募集】多様なプロジェクトに参画可能!
$ws.Range("A17").Value = '2025-12-17 01:20:50'
$ws.Range("B17").Value = '【SESエンジニア募集】多様なプロジェクトに参画可能!'
$ws.Range("C17").Value = 'システム開発'
$ws.Range("D17").Value = '300,000 円 ~ 500,000 円 / 固定'
$ws.Range("E17").Value = '期限情報なし'
$ws.Range("F17").Value = 'https://www.lancers.jp/work/detail/5437544'
$ws.Range("G17").Value = 25
$ws.Range("H17").Value = ""

# Recreate hyperlinks for column F (URL) so each cell links to its own displayed URL
$ws.Hyperlinks.Add($ws.Range("F2"), 'https://www.lancers.jp/work/detail/5455098')
$ws.Hyperlinks.Add($ws.Range("F3"), 'https://www.lancers.jp/work/detail/5454985')
$ws.Hyperlinks.Add($ws.Range("F4"), 'https://www.lancers.jp/work/detail/5217096')
$ws.Hyperlinks.Add($ws.Range("F5"), 'https://www.lancers.jp/work/detail/5439921')
$ws.Hyperlinks.Add($ws.Range("F6"), 'https://www.lancers.jp/work/detail/5455160')
$ws.Hyperlinks.Add($ws.Range("F7"), 'https://www.lancers.jp/work/detail/5455415')
$ws.Hyperlinks.Add($ws.Range("F8"), 'https://www.lancers.jp/work/detail/5455038')
$ws.Hyperlinks.Add($ws.Range("F9"), 'https://www.lancers.jp/work/detail/5455251')
$ws.Hyperlinks.Add($ws.Range("F10"), 'https://www.lancers.jp/work/detail/5455029')
$ws.Hyperlinks.Add($ws.Range("F11"), 'https://www.lancers.jp/work/detail/5454857')
$ws.Hyperlinks.Add($ws.Range("F12"), 'https://www.lancers.jp/work/detail/5455015')
$ws.Hyperlinks.Add($ws.Range("F13"), 'https://www.lancers.jp/work/detail/5455422')
$ws.Hyperlinks.Add($ws.Range("F14"), 'https://www.lancers.jp/work/detail/5455067')
$ws.Hyperlinks.Add($ws.Range("F15"), 'https://www.lancers.jp/work/detail/5450345')
$ws.Hyperlinks.Add($ws.Range("F16"), 'https://www.lancers.jp/work/detail/5341051')
$ws.Hyperlinks.Add($ws.Range("F17"), 'https://www.lancers.jp/work/detail/5437544')

# Re-apply the Hyperlink cell style (Hyperlinks.Add() creates a duplicate style entry otherwise)
$ws.Range("F2").Style = "Hyperlink"
$ws.Range("F3").Style = "Hyperlink"
$ws.Range("F4").Style = "Hyperlink"
$ws.Range("F5").Style = "Hyperlink"
$ws.Range("F6").Style = "Hyperlink"
$ws.Range("F7").Style = "Hyperlink"
$ws.Range("F8").Style = "Hyperlink"
$ws.Range("F9").Style = "Hyperlink"
$ws.Range("F10").Style = "Hyperlink"
$ws.Range("F11").Style = "Hyperlink"
$ws.Range("F12").Style = "Hyperlink"
$ws.Range("F13").Style = "Hyperlink"
$ws.Range("F14").Style = "Hyperlink"
$ws.Range("F15").Style = "Hyperlink"
$ws.Range("F16").Style = "Hyperlink"
$ws.Range("F17").Style = "Hyperlink"
